$wb = $excel.ActiveWorkbook

# Grab a reference to the existing sheet before mutating the sheets collection.
$odiSheet = $wb.Worksheets.Item("ODI Batting")

# Add a brand new worksheet; Add() with no args inserts it before the active
# sheet, i.e. at position 1 - exactly where "Player Info" needs to land.
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Sheet handles returned by the COM layer are positional, so inserting a new
# sheet ahead of "ODI Batting" invalidates the earlier $odiSheet handle
# (it now points at whatever sheet occupies that old slot). Re-resolve both
# sheets by name now that the collection is stable.
$odiSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Item("Player Info")

# ---- Populate the new "Player Info" sheet ----
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold / centered / bordered header look used on "ODI Batting".
# (A cross-sheet `.Style =` assignment doesn't stick, so copy formats only.)
$odiSheet.Range("A1:D1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# ID is textual data ("5500"), not a number - lead with an apostrophe so
# Excel stores it as text, then drop back to the Normal style so the
# quote-prefix flag doesn't linger as a visible cell format.
$playerInfo.Range("A2").Value = "'5500"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Minod Bhanuka Ranasinghe"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

# ---- "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, full URL -> bare code ----
$odiSheet.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{ 2 = "4376"; 3 = "4480"; 4 = "4482"; 5 = "4485"; 6 = "4487"; 7 = "4488" }
foreach ($row in $matchCodes.Keys) {
    $cell = $odiSheet.Range("D" + $row)
    $cell.Value = "'" + $matchCodes[$row]
    $cell.Style = "Normal"
}
